# Update TPM-derived NATMI metrics on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ G = 3.627806666666667;  H = 10.88342;          I = 0.08036675778141429; J = 0.08036675778141429;
            M = 2.724001666666667;  N = 8.172005;           O = 0.04635500474236593; P = 0.04635500474236593;
            Q = 9.882151406344445;  R = 88.93936265710001;  S = 0.003725401438086034; T = 0.003725401438086034 }
    3  = @{ G = 3.627806666666667;  H = 10.88342;          I = 0.08036675778141429; J = 0.08036675778141429;
            O = 0.6912512390256352; P = 0.6912512390256351;
            Q = 147.3637947367378;  R = 1326.27415263064;   S = 0.05555362089287574;  T = 0.05555362089287573 }
    4  = @{ G = 3.627806666666667;  H = 10.88342;          I = 0.08036675778141429; J = 0.08036675778141429;
            M = 15.419285;          N = 46.257855;          O = 0.2623937562319988;  P = 0.2623937562319988;
            Q = 55.93818491823334;  R = 503.4436642641001;  S = 0.02108773545045252;  T = 0.02108773545045252 }
    5  = @{ I = 0.6506403335968259; J = 0.6506403335968259;
            M = 2.724001666666667;  N = 8.172005;           O = 0.04635500474236593; P = 0.04635500474236593;
            Q = 80.00479881453222;  R = 720.04318933079;    S = 0.03016043574945542;  T = 0.03016043574945542 }
    6  = @{ I = 0.6506403335968259; J = 0.6506403335968259;
            O = 0.6912512390256352; P = 0.6912512390256351;
            S = 0.4497559367588586; T = 0.4497559367588585 }
    7  = @{ I = 0.6506403335968259; J = 0.6506403335968259;
            M = 15.419285;          N = 46.257855;          O = 0.2623937562319988;  P = 0.2623937562319988;
            Q = 452.8693243416766;  R = 4075.82391907509;   S = 0.170723961088512;    T = 0.170723961088512 }
    8  = @{ G = 12.14251133333333; H = 36.427534;          I = 0.2689929086217598;  J = 0.2689929086217598;
            M = 2.724001666666667;  N = 8.172005;           O = 0.04635500474236593; P = 0.04635500474236593;
            Q = 33.07622110951889;  R = 297.68598998567;    S = 0.01246916755482448;  T = 0.01246916755482448 }
    9  = @{ G = 12.14251133333333; H = 36.427534;          I = 0.2689929086217598;  J = 0.2689929086217598;
            O = 0.6912512390256352; P = 0.6912512390256351;
            Q = 493.2364682371476;  R = 4439.128214134329;  S = 0.185941681373901;    T = 0.1859416813739009 }
    10 = @{ G = 12.14251133333333; H = 36.427534;          I = 0.2689929086217598;  J = 0.2689929086217598;
            M = 15.419285;          N = 46.257855;          O = 0.2623937562319988;  P = 0.2623937562319988;
            Q = 187.2288428643967;  R = 1685.05958577957;   S = 0.0705820596930344;   T = 0.0705820596930344 }
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
